# Apply inventory issuance transaction:
#  - Update "Current Stock" sheet quantities/metadata for White Tiles and Cement
#  - Append two new rows to the "Issuance Log" sheet recording the issuance

$wb = $excel.ActiveWorkbook

$stock = $wb.Worksheets.Item("Current Stock")
$log = $wb.Worksheets.Item("Issuance Log")

# --- Update Current Stock sheet ---
# Row 2: White Tiles
$stock.Range("B2").Value = 0.5
$stock.Range("D2").Value = "2025-07-10 12:45:18"
$stock.Range("E2").Value = "engineer"

# Row 3: Cement (50 Kgs)
$stock.Range("B3").Value = 80
$stock.Range("D3").Value = "2025-07-10 12:45:18"
$stock.Range("E3").Value = "engineer"

# --- Append new rows to Issuance Log sheet ---
# Row 3: issuance of White Tiles
$log.Range("A3").Value = "2025-07-10 12:45:18"
$log.Range("B3").Value = "White Tiles"
$log.Range("C3").Value = 0.5
$log.Range("D3").Value = "boxes"
$log.Range("E3").Value = "engineer"
$log.Range("F3").Value = "Batch BTH-20250710-624B - Authorized by engineer"
$log.Range("G3").Value = 0.5

# Row 4: issuance of Cement (50 Kgs)
$log.Range("A4").Value = "2025-07-10 12:45:18"
$log.Range("B4").Value = "Cement (50 Kgs)"
$log.Range("C4").Value = 20
$log.Range("D4").Value = "bags"
$log.Range("E4").Value = "engineer"
$log.Range("F4").Value = "Batch BTH-20250710-624B - Authorized by engineer"
$log.Range("G4").Value = 80
